$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 "Marking": Right marks and Wrong marks corrected
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total": recompute total marks and displayed fraction to fix total marks error
$ws.Range("B12").Value = 84
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "82 / 112"
